$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value2 = 46056.01041666666
$ws.Cells.Item(2, 2).Value2 = 154.003
$ws.Cells.Item(3, 1).Value2 = 46056.02083333334
$ws.Cells.Item(3, 2).Value2 = 148.826
$ws.Cells.Item(4, 1).Value2 = 46056.03125
$ws.Cells.Item(4, 2).Value2 = 144.383
$ws.Cells.Item(5, 1).Value2 = 46056.04166666666
$ws.Cells.Item(5, 2).Value2 = 139.561
$ws.Cells.Item(6, 1).Value2 = 46056.05208333334
$ws.Cells.Item(6, 2).Value2 = 131.369
$ws.Cells.Item(7, 1).Value2 = 46056.0625
$ws.Cells.Item(7, 2).Value2 = 126.962
$ws.Cells.Item(8, 1).Value2 = 46056.07291666666
$ws.Cells.Item(8, 2).Value2 = 122.441
$ws.Cells.Item(9, 1).Value2 = 46056.08333333334
$ws.Cells.Item(9, 2).Value2 = 118.598
$ws.Cells.Item(10, 1).Value2 = 46056.09375
$ws.Cells.Item(10, 2).Value2 = 114.35
$ws.Cells.Item(11, 1).Value2 = 46056.10416666666
$ws.Cells.Item(11, 2).Value2 = 111.643
$ws.Cells.Item(12, 1).Value2 = 46056.11458333334
$ws.Cells.Item(12, 2).Value2 = 109.267
$ws.Cells.Item(13, 1).Value2 = 46056.125
$ws.Cells.Item(13, 2).Value2 = 107.022
$ws.Cells.Item(14, 1).Value2 = 46056.13541666666
$ws.Cells.Item(14, 2).Value2 = 103.261
$ws.Cells.Item(15, 1).Value2 = 46056.14583333334
$ws.Cells.Item(15, 2).Value2 = 101.261
$ws.Cells.Item(16, 1).Value2 = 46056.15625
$ws.Cells.Item(16, 2).Value2 = 98.959
$ws.Cells.Item(17, 1).Value2 = 46056.16666666666
$ws.Cells.Item(17, 2).Value2 = 97.19499999999999
$ws.Cells.Item(18, 1).Value2 = 46056.17708333334
$ws.Cells.Item(18, 2).Value2 = 73.86
$ws.Cells.Item(19, 1).Value2 = 46056.1875
$ws.Cells.Item(19, 2).Value2 = 104.835
$ws.Cells.Item(20, 1).Value2 = 46056.19791666666
$ws.Cells.Item(20, 2).Value2 = 103.351
$ws.Cells.Item(21, 1).Value2 = 46056.20833333334
$ws.Cells.Item(21, 2).Value2 = 102.08
$ws.Cells.Item(22, 1).Value2 = 46056.21875
$ws.Cells.Item(22, 2).Value2 = 77.923
$ws.Cells.Item(23, 1).Value2 = 46056.22916666666
$ws.Cells.Item(23, 2).Value2 = 76.53100000000001
$ws.Cells.Item(24, 1).Value2 = 46056.23958333334
$ws.Cells.Item(24, 2).Value2 = 75.20099999999999
$ws.Cells.Item(25, 1).Value2 = 46056.25
$ws.Cells.Item(25, 2).Value2 = 107.297
$ws.Cells.Item(26, 1).Value2 = 46056.26041666666
$ws.Cells.Item(26, 2).Value2 = 76.929
$ws.Cells.Item(27, 1).Value2 = 46056.27083333334
$ws.Cells.Item(27, 2).Value2 = 76.789
$ws.Cells.Item(28, 1).Value2 = 46056.28125
$ws.Cells.Item(28, 2).Value2 = 76.97499999999999
$ws.Cells.Item(29, 1).Value2 = 46056.29166666666
$ws.Cells.Item(29, 2).Value2 = 77.18600000000001
$ws.Cells.Item(30, 1).Value2 = 46056.30208333334
$ws.Cells.Item(30, 2).Value2 = 69.506
$ws.Cells.Item(31, 1).Value2 = 46056.3125
$ws.Cells.Item(31, 2).Value2 = 70.932
$ws.Cells.Item(32, 1).Value2 = 46056.32291666666
$ws.Cells.Item(32, 2).Value2 = 72.7
$ws.Cells.Item(33, 1).Value2 = 46056.33333333334
$ws.Cells.Item(33, 2).Value2 = 81.31999999999999
$ws.Cells.Item(34, 1).Value2 = 46056.34375
$ws.Cells.Item(34, 2).Value2 = 73.096
$ws.Cells.Item(35, 1).Value2 = 46056.35416666666
$ws.Cells.Item(35, 2).Value2 = 74.39400000000001
$ws.Cells.Item(36, 1).Value2 = 46056.36458333334
$ws.Cells.Item(36, 2).Value2 = 77.167
$ws.Cells.Item(37, 1).Value2 = 46056.375
$ws.Cells.Item(37, 2).Value2 = 77.78400000000001
$ws.Cells.Item(38, 1).Value2 = 46056.38541666666
$ws.Cells.Item(38, 2).Value2 = 104.846
$ws.Cells.Item(39, 1).Value2 = 46056.39583333334
$ws.Cells.Item(39, 2).Value2 = 105.444
$ws.Cells.Item(40, 1).Value2 = 46056.40625
$ws.Cells.Item(40, 2).Value2 = 105.095
$ws.Cells.Item(41, 1).Value2 = 46056.41666666666
$ws.Cells.Item(41, 2).Value2 = 105.826
$ws.Cells.Item(42, 1).Value2 = 46056.42708333334
$ws.Cells.Item(42, 2).Value2 = 111.221
$ws.Cells.Item(43, 1).Value2 = 46056.4375
$ws.Cells.Item(43, 2).Value2 = 112.251
$ws.Cells.Item(44, 1).Value2 = 46056.44791666666
$ws.Cells.Item(44, 2).Value2 = 113.463
$ws.Cells.Item(45, 1).Value2 = 46056.45833333334
$ws.Cells.Item(45, 2).Value2 = 114.303
$ws.Cells.Item(46, 1).Value2 = 46056.46875
$ws.Cells.Item(46, 2).Value2 = 115.169
$ws.Cells.Item(47, 1).Value2 = 46056.47916666666
$ws.Cells.Item(47, 2).Value2 = 116.636
$ws.Cells.Item(48, 1).Value2 = 46056.48958333334
$ws.Cells.Item(48, 2).Value2 = 117.575
$ws.Cells.Item(49, 1).Value2 = 46056.5
$ws.Cells.Item(49, 2).Value2 = 118.651
$ws.Cells.Item(50, 1).Value2 = 46056.51041666666
$ws.Cells.Item(50, 2).Value2 = 122.039
$ws.Cells.Item(51, 1).Value2 = 46056.52083333334
$ws.Cells.Item(51, 2).Value2 = 123.929
$ws.Cells.Item(52, 1).Value2 = 46056.53125
$ws.Cells.Item(52, 2).Value2 = 126.78
$ws.Cells.Item(53, 1).Value2 = 46056.54166666666
$ws.Cells.Item(53, 2).Value2 = 128.719
$ws.Cells.Item(54, 1).Value2 = 46056.55208333334
$ws.Cells.Item(54, 2).Value2 = 134.09
$ws.Cells.Item(55, 1).Value2 = 46056.5625
$ws.Cells.Item(55, 2).Value2 = 137.724
$ws.Cells.Item(56, 1).Value2 = 46056.57291666666
$ws.Cells.Item(56, 2).Value2 = 142.514
$ws.Cells.Item(57, 1).Value2 = 46056.58333333334
$ws.Cells.Item(57, 2).Value2 = 145.959
$ws.Cells.Item(58, 1).Value2 = 46056.59375
$ws.Cells.Item(58, 2).Value2 = 152.466
$ws.Cells.Item(59, 1).Value2 = 46056.60416666666
$ws.Cells.Item(59, 2).Value2 = 157.746
$ws.Cells.Item(60, 1).Value2 = 46056.61458333334
$ws.Cells.Item(60, 2).Value2 = 164.443
$ws.Cells.Item(61, 1).Value2 = 46056.625
$ws.Cells.Item(61, 2).Value2 = 169.527
$ws.Cells.Item(62, 1).Value2 = 46056.63541666666
$ws.Cells.Item(62, 2).Value2 = 184.667
$ws.Cells.Item(63, 1).Value2 = 46056.64583333334
$ws.Cells.Item(63, 2).Value2 = 196.697
$ws.Cells.Item(64, 1).Value2 = 46056.65625
$ws.Cells.Item(64, 2).Value2 = 208.533
$ws.Cells.Item(65, 1).Value2 = 46056.66666666666
$ws.Cells.Item(65, 2).Value2 = 246.418
$ws.Cells.Item(66, 1).Value2 = 46056.67708333334
$ws.Cells.Item(66, 2).Value2 = 253.584
$ws.Cells.Item(67, 1).Value2 = 46056.6875
$ws.Cells.Item(67, 2).Value2 = 281.775
$ws.Cells.Item(68, 1).Value2 = 46056.69791666666
$ws.Cells.Item(68, 2).Value2 = 311.001
$ws.Cells.Item(69, 1).Value2 = 46056.70833333334
$ws.Cells.Item(69, 2).Value2 = 339.528
$ws.Cells.Item(70, 1).Value2 = 46056.71875
$ws.Cells.Item(70, 2).Value2 = 387.296
$ws.Cells.Item(71, 1).Value2 = 46056.72916666666
$ws.Cells.Item(71, 2).Value2 = 418.698
$ws.Cells.Item(72, 1).Value2 = 46056.73958333334
$ws.Cells.Item(72, 2).Value2 = 450.255
$ws.Cells.Item(73, 1).Value2 = 46056.75
$ws.Cells.Item(73, 2).Value2 = 481.633
$ws.Cells.Item(74, 1).Value2 = 46056.76041666666
$ws.Cells.Item(74, 2).Value2 = 530.521
$ws.Cells.Item(75, 1).Value2 = 46056.77083333334
$ws.Cells.Item(75, 2).Value2 = 558.141
$ws.Cells.Item(76, 1).Value2 = 46056.78125
$ws.Cells.Item(76, 2).Value2 = 586.798
$ws.Cells.Item(77, 1).Value2 = 46056.79166666666
$ws.Cells.Item(77, 2).Value2 = 614.644
$ws.Cells.Item(78, 1).Value2 = 46056.80208333334
$ws.Cells.Item(78, 2).Value2 = 653.846
$ws.Cells.Item(79, 1).Value2 = 46056.8125
$ws.Cells.Item(79, 2).Value2 = 672.149
$ws.Cells.Item(80, 1).Value2 = 46056.82291666666
$ws.Cells.Item(80, 2).Value2 = 689.759
$ws.Cells.Item(81, 1).Value2 = 46056.83333333334
$ws.Cells.Item(81, 2).Value2 = 707.6559999999999
$ws.Cells.Item(82, 1).Value2 = 46056.84375
$ws.Cells.Item(82, 2).Value2 = 741.17
$ws.Cells.Item(83, 1).Value2 = 46056.85416666666
$ws.Cells.Item(83, 2).Value2 = 757.683
$ws.Cells.Item(84, 1).Value2 = 46056.86458333334
$ws.Cells.Item(84, 2).Value2 = 774.265
$ws.Cells.Item(85, 1).Value2 = 46056.875
$ws.Cells.Item(85, 2).Value2 = 791.246
$ws.Cells.Item(86, 1).Value2 = 46056.88541666666
$ws.Cells.Item(86, 2).Value2 = 809.602
$ws.Cells.Item(87, 1).Value2 = 46056.89583333334
$ws.Cells.Item(87, 2).Value2 = 821.9880000000001
$ws.Cells.Item(88, 1).Value2 = 46056.90625
$ws.Cells.Item(88, 2).Value2 = 834.153
$ws.Cells.Item(89, 1).Value2 = 46056.91666666666
$ws.Cells.Item(89, 2).Value2 = 846.282
$ws.Cells.Item(90, 1).Value2 = 46056.92708333334
$ws.Cells.Item(90, 2).Value2 = 862.717
$ws.Cells.Item(91, 1).Value2 = 46056.9375
$ws.Cells.Item(91, 2).Value2 = 875.299
$ws.Cells.Item(92, 1).Value2 = 46056.94791666666
$ws.Cells.Item(92, 2).Value2 = 888.052
$ws.Cells.Item(93, 1).Value2 = 46056.95833333334
$ws.Cells.Item(93, 2).Value2 = 900.475
$ws.Cells.Item(94, 1).Value2 = 46056.96875
$ws.Cells.Item(94, 2).Value2 = 0
$ws.Cells.Item(95, 1).Value2 = 46056.97916666666
$ws.Cells.Item(95, 2).Value2 = 0
$ws.Cells.Item(96, 1).Value2 = 46056.98958333334
$ws.Cells.Item(96, 2).Value2 = 0
$ws.Cells.Item(97, 1).Value2 = 46057
$ws.Cells.Item(97, 2).Value2 = 0
